{"js": "// Replace the date line and every \"NNN\u00d7N=\" multiplication expression in the\n// document with the updated values described by the commit diff.\nconst replacements = [\n  [\"2026-01-18 Sunday\", \"2026-01-19 Monday\"],\n  [\"494\u00d79=\", \"305\u00d77=\"],\n  [\"335\u00d74=\", \"114\u00d76=\"],\n  [\"734\u00d77=\", \"264\u00d78=\"],\n  [\"179\u00d75=\", \"816\u00d74=\"],\n  [\"190\u00d76=\", \"312\u00d73=\"],\n  [\"428\u00d78=\", \"482\u00d73=\"],\n  [\"252\u00d75=\", \"944\u00d74=\"],\n  [\"514\u00d79=\", \"586\u00d77=\"],\n  [\"901\u00d77=\", \"464\u00d74=\"],\n  [\"311\u00d72=\", \"907\u00d76=\"],\n  [\"772\u00d76=\", \"951\u00d75=\"],\n  [\"686\u00d78=\", \"428\u00d79=\"],\n  [\"622\u00d76=\", \"304\u00d74=\"],\n  [\"444\u00d78=\", \"710\u00d78=\"],\n  [\"760\u00d76=\", \"180\u00d75=\"],\n  [\"930\u00d77=\", \"129\u00d72=\"],\n  [\"317\u00d75=\", \"828\u00d79=\"],\n  [\"982\u00d75=\", \"596\u00d73=\"],\n  [\"863\u00d78=\", \"250\u00d72=\"],\n  [\"806\u00d77=\", \"859\u00d76=\"],\n  [\"431\u00d79=\", \"493\u00d73=\"],\n  [\"546\u00d78=\", \"928\u00d75=\"],\n  [\"480\u00d74=\", \"190\u00d77=\"],\n  [\"714\u00d76=\", \"795\u00d77=\"],\n  [\"201\u00d77=\", \"836\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2026-01-18 Sunday\", \"2026-01-19 Monday\"),\n  @(\"494\u00d79=\", \"305\u00d77=\"),\n  @(\"335\u00d74=\", \"114\u00d76=\"),\n  @(\"734\u00d77=\", \"264\u00d78=\"),\n  @(\"179\u00d75=\", \"816\u00d74=\"),\n  @(\"190\u00d76=\", \"312\u00d73=\"),\n  @(\"428\u00d78=\", \"482\u00d73=\"),\n  @(\"252\u00d75=\", \"944\u00d74=\"),\n  @(\"514\u00d79=\", \"586\u00d77=\"),\n  @(\"901\u00d77=\", \"464\u00d74=\"),\n  @(\"311\u00d72=\", \"907\u00d76=\"),\n  @(\"772\u00d76=\", \"951\u00d75=\"),\n  @(\"686\u00d78=\", \"428\u00d79=\"),\n  @(\"622\u00d76=\", \"304\u00d74=\"),\n  @(\"444\u00d78=\", \"710\u00d78=\"),\n  @(\"760\u00d76=\", \"180\u00d75=\"),\n  @(\"930\u00d77=\", \"129\u00d72=\"),\n  @(\"317\u00d75=\", \"828\u00d79=\"),\n  @(\"982\u00d75=\", \"596\u00d73=\"),\n  @(\"863\u00d78=\", \"250\u00d72=\"),\n  @(\"806\u00d77=\", \"859\u00d76=\"),\n  @(\"431\u00d79=\", \"493\u00d73=\"),\n  @(\"546\u00d78=\", \"928\u00d75=\"),\n  @(\"480\u00d74=\", \"190\u00d77=\"),\n  @(\"714\u00d76=\", \"795\u00d77=\"),\n  @(\"201\u00d77=\", \"836\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\n$d.Save()\n"}
